$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text-valued cells (names, URLs, labels) - safe to assign directly
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E21").Value = "20HotbitTokenHTBBestin24h"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

# Numeric-looking but text-semantic cells (Price/Hora columns stored as text in source)
# Force text format so Excel does not auto-convert to a Number type, then restore default style
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "246.18"
$c.Style = "Normal"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "24.15"
$c.Style = "Normal"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.295"
$c.Style = "Normal"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.05760"
$c.Style = "Normal"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "6.476"
$c.Style = "Normal"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.138"
$c.Style = "Normal"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.8191"
$c.Style = "Normal"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.8803"
$c.Style = "Normal"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.1378"
$c.Style = "Normal"
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.06957"
$c.Style = "Normal"
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.03134"
$c.Style = "Normal"
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.02938"
$c.Style = "Normal"
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.09392"
$c.Style = "Normal"
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.746"
$c.Style = "Normal"
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.001529"
$c.Style = "Normal"
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.04716"
$c.Style = "Normal"
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.0006010"
$c.Style = "Normal"
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.006232"
$c.Style = "Normal"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.001242"
$c.Style = "Normal"
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.004664"
$c.Style = "Normal"
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.00008798"
$c.Style = "Normal"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.531"
$c.Style = "Normal"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.152"
$c.Style = "Normal"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.1312"
$c.Style = "Normal"
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.03720"
$c.Style = "Normal"
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1057"
$c.Style = "Normal"
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.002739"
$c.Style = "Normal"
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.003072"
$c.Style = "Normal"
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.008299"
$c.Style = "Normal"
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00005243"
$c.Style = "Normal"
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.002214"
$c.Style = "Normal"
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = "16"
$c.Style = "Normal"
